$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header C1 text: crop yield -> crop production index
$ws.Range("C1").Value = "0. Crop production index"

# Add new header I1 (copy H1's formatting/style, then set its text)
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "17. Employment in agriculture (% of total employment) (modeled ILO estimate)"

# Row 2 becomes the old row 4 (2014) data; row 2's crop value is now a text string, not a number
$ws.Range("B2").Value = 2014
$ws.Range("D2").Value = 27.65
$ws.Range("E2").Value = 59.66004283
$ws.Range("F2").Value = 235.09
$ws.Range("G2").Value = 4.200625935
$ws.Range("H2").Value = 37003245
$ws.Range("I2").Value = 43.2859542809493

# Force C2 to hold the text string "106.91" (not an auto-converted number) while
# keeping the cell's original (unformatted) style, by writing a text formula and
# then collapsing it down to its static value via copy/paste-values.
$ws.Range("C2").Formula = '=T("106.91")'
$ws.Range("C2").Copy()
$ws.Range("C2").PasteSpecial(-4163)  # xlPasteValues

# Drop the old rows 3 and 4 (2013 and 2014 data that is no longer needed,
# since row 2 now carries the 2014 figures and the sheet only spans 2 rows)
$ws.Range("A3:I4").EntireRow.Delete()
